$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 17 de Junio de 2020 a las 02:08"

# --- Reorder country rows (shared-string content swap) ---
$ws.Cells.Item(91,1).Value = "Venezuela"
$ws.Cells.Item(92,1).Value = "Grecia"
$ws.Cells.Item(93,1).Value = "Tailandia"
$ws.Cells.Item(94,1).Value = "Bosnia y Herzegovina"
$ws.Cells.Item(206,1).Value = "Groenlandia"
$ws.Cells.Item(207,1).Value = "Islas Malvinas"
$ws.Cells.Item(210,1).Value = "Seychelles"
$ws.Cells.Item(211,1).Value = "Montserrat"
$ws.Cells.Item(213,1).Value = "Papua Nueva Guinea"
$ws.Cells.Item(214,1).Value = "Islas Virgenes Britanicas"

# --- Update daily statistics ---
# Row 4
$ws.Cells.Item(4,2).Value = 2208241
$ws.Cells.Item(4,3).Value = 25291
$ws.Cells.Item(4,4).Value = 899503
$ws.Cells.Item(4,5).Value = 1189609
$ws.Cells.Item(4,7).Value = 846
$ws.Cells.Item(4,8).Value = 119129
# Row 5
$ws.Cells.Item(5,2).Value = 928834
$ws.Cells.Item(5,3).Value = 37278
$ws.Cells.Item(5,5).Value = 418604
$ws.Cells.Item(5,7).Value = 1338
$ws.Cells.Item(5,8).Value = 45456
# Row 20
$ws.Cells.Item(20,2).Value = 99467
$ws.Cells.Item(20,3).Value = 320
$ws.Cells.Item(20,4).Value = 61443
$ws.Cells.Item(20,5).Value = 29811
# Row 37
$ws.Cells.Item(37,2).Value = 34159
$ws.Cells.Item(37,3).Value = 1374
$ws.Cells.Item(37,5).Value = 23107
$ws.Cells.Item(37,7).Value = 24
$ws.Cells.Item(37,8).Value = 878
# Row 48
$ws.Cells.Item(48,2).Value = 21962
$ws.Cells.Item(48,3).Value = 540
$ws.Cells.Item(48,5).Value = 7739
$ws.Cells.Item(48,7).Value = 9
$ws.Cells.Item(48,8).Value = 457
# Row 55
$ws.Cells.Item(55,2).Value = 17148
$ws.Cells.Item(55,3).Value = 490
$ws.Cells.Item(55,4).Value = 5623
$ws.Cells.Item(55,5).Value = 11070
$ws.Cells.Item(55,7).Value = 31
$ws.Cells.Item(55,8).Value = 455
# Row 91
$ws.Cells.Item(91,2).Value = 3150
$ws.Cells.Item(91,3).Value = 88
$ws.Cells.Item(91,4).Value = 835
$ws.Cells.Item(91,5).Value = 2288
$ws.Cells.Item(91,8).Value = 27
# Row 92
$ws.Cells.Item(92,2).Value = 3148
$ws.Cells.Item(92,3).Value = 14
$ws.Cells.Item(92,4).Value = 1374
$ws.Cells.Item(92,5).Value = 1589
$ws.Cells.Item(92,7).Value = 1
$ws.Cells.Item(92,8).Value = 185
# Row 93
$ws.Cells.Item(93,2).Value = 3135
$ws.Cells.Item(93,3).Value = 0
$ws.Cells.Item(93,4).Value = 2993
$ws.Cells.Item(93,5).Value = 84
$ws.Cells.Item(93,7).Value = 0
$ws.Cells.Item(93,8).Value = 58
# Row 94
$ws.Cells.Item(94,2).Value = 3085
$ws.Cells.Item(94,3).Value = 45
$ws.Cells.Item(94,4).Value = 2178
$ws.Cells.Item(94,5).Value = 739
$ws.Cells.Item(94,7).Value = 3
$ws.Cells.Item(94,8).Value = 168
# Row 135
$ws.Cells.Item(135,2).Value = 849
$ws.Cells.Item(135,3).Value = 1
$ws.Cells.Item(135,4).Value = 801
$ws.Cells.Item(135,5).Value = 24
$ws.Cells.Item(135,7).Value = 1
$ws.Cells.Item(135,8).Value = 24
# Row 152
$ws.Cells.Item(152,2).Value = 495
$ws.Cells.Item(152,5).Value = 34
# Row 153
$ws.Cells.Item(153,2).Value = 484
$ws.Cells.Item(153,3).Value = 17
$ws.Cells.Item(153,4).Value = 76
$ws.Cells.Item(153,5).Value = 398
# Row 210
$ws.Cells.Item(210,4).Value = 11
$ws.Cells.Item(210,8).Value = 0
# Row 211
$ws.Cells.Item(211,4).Value = 10
$ws.Cells.Item(211,8).Value = 1
# Row 213
$ws.Cells.Item(213,4).Value = 8
$ws.Cells.Item(213,8).Value = 0
# Row 214
$ws.Cells.Item(214,4).Value = 7
$ws.Cells.Item(214,8).Value = 1
